# admissionCreatePatient.xlsx - "update for scheduling testcase"
#
# Adds a new "Default Value" column (H) to the testcase sheet and fills in
# the default values used for the First Name / Last Name steps, then
# leaves the selection on the newly-edited cell (H7), matching the
# scheduling testcase update captured in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for column H.
$ws.Range("H1").Value = "Default Value"

# Default values tied to the "First Name" (row 6) and "Last Name" (row 7)
# input steps.
$ws.Range("H6").Value = "firstname"
$ws.Range("H7").Value = "lastname"

# Row 6 grows taller (15 -> 30) to accommodate the new content.
$ws.Rows.Item(6).RowHeight = 30

# Leave the selection on H7, as in the updated workbook.
$ws.Range("H7").Select()
